$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 27028702
$ws.Range("I137").Value = 1120.2142
$ws.Range("K137").Value = 3360.6426
$ws.Range("M137").Value = -810.6425999999997

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2452080.5
$ws.Range("I2").Value = 500
$ws.Range("J2").Value = 4203209.5
$ws.Range("K2").Value = 500
$ws.Range("L2").Value = 4203209.5
$ws.Range("M2").Value = -387
$ws.Range("N2").Value = -4203435.5

$ws.Range("H45").Value = 749.5789
$ws.Range("I45").Value = 637.4545000000001
$ws.Range("J45").Value = 903.75
$ws.Range("K45").Value = 637.4545000000001
$ws.Range("L45").Value = 903.75
$ws.Range("M45").Value = -260.4545000000001
$ws.Range("N45").Value = -1657.75

$ws.Range("H61").Value = 3686.8
$ws.Range("I61").Value = 2156.6667
$ws.Range("J61").Value = 5982
$ws.Range("K61").Value = 2156.6667
$ws.Range("L61").Value = 5982
$ws.Range("M61").Value = -1944.6667
$ws.Range("N61").Value = -6406

$ws.Range("H74").Value = 5688.3706
$ws.Range("I74").Value = 1368.6666
$ws.Range("J74").Value = 6922.5713
$ws.Range("K74").Value = 1368.6666
$ws.Range("L74").Value = 6922.5713
$ws.Range("M74").Value = -494.6666
$ws.Range("N74").Value = -8670.5713

$ws.Range("H77").Value = 5688.3706
$ws.Range("I77").Value = 1368.6666
$ws.Range("J77").Value = 6922.5713
$ws.Range("K77").Value = 6843.333000000001
$ws.Range("L77").Value = 34612.85649999999
$ws.Range("M77").Value = -2475.333000000001
$ws.Range("N77").Value = -43348.85649999999

$ws.Range("H97").Value = 545.625
$ws.Range("I97").Value = 493.76923
$ws.Range("J97").Value = 770.3333
$ws.Range("K97").Value = 493.76923
$ws.Range("L97").Value = 770.3333
$ws.Range("M97").Value = 2.230770000000007
$ws.Range("N97").Value = -1762.3333

$ws.Range("H110").Value = 1436
$ws.Range("I110").Value = 1236
$ws.Range("K110").Value = 1236
$ws.Range("M110").Value = 809

$ws.Range("H116").Value = 2452080.5
$ws.Range("I116").Value = 500
$ws.Range("J116").Value = 4203209.5
$ws.Range("K116").Value = 500
$ws.Range("L116").Value = 4203209.5
$ws.Range("M116").Value = 1794
$ws.Range("N116").Value = -4207797.5

$ws.Range("H122").Value = 1940.7273
$ws.Range("I122").Value = 1692.5714
$ws.Range("J122").Value = 2375
$ws.Range("K122").Value = 5077.7142
$ws.Range("L122").Value = 7125
$ws.Range("M122").Value = -2627.7142
$ws.Range("N122").Value = -12025

$ws.Range("H132").Value = 1562.8235
$ws.Range("I132").Value = 1342.6897
$ws.Range("J132").Value = 2839.6
$ws.Range("K132").Value = 4028.0691
$ws.Range("L132").Value = 8518.799999999999
$ws.Range("M132").Value = -1498.0691
$ws.Range("N132").Value = -13578.8

$ws.Range("H136").Value = 3686.8
$ws.Range("I136").Value = 2156.6667
$ws.Range("J136").Value = 5982
$ws.Range("K136").Value = 6470.000100000001
$ws.Range("L136").Value = 17946
$ws.Range("M136").Value = -3920.000100000001
$ws.Range("N136").Value = -23046

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2452080.5
$ws.Range("I3").Value = 500
$ws.Range("J3").Value = 4203209.5
$ws.Range("K3").Value = 500
$ws.Range("L3").Value = 4203209.5
$ws.Range("M3").Value = -386
$ws.Range("N3").Value = -4203437.5

$ws.Range("H99").Value = 1807.1428
$ws.Range("I99").Value = 1630
$ws.Range("K99").Value = 1630
$ws.Range("M99").Value = -132

$ws.Range("H107").Value = 1760.8948
$ws.Range("I107").Value = 1675.3889
$ws.Range("J107").Value = 3300
$ws.Range("K107").Value = 1675.3889
$ws.Range("L107").Value = 3300
$ws.Range("M107").Value = 244.6111000000001
$ws.Range("N107").Value = -7140

$ws.Range("H134").Value = 40331.656
$ws.Range("I134").Value = 42763.11
$ws.Range("J134").Value = 7507
$ws.Range("K134").Value = 128289.33
$ws.Range("L134").Value = 22521
$ws.Range("M134").Value = -125754.33
$ws.Range("N134").Value = -27591

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1214.2667
$ws.Range("I31").Value = 1054.5834
$ws.Range("K31").Value = 1054.5834
$ws.Range("M31").Value = -759.5834

$ws.Range("H34").Value = 1214.2667
$ws.Range("I34").Value = 1054.5834
$ws.Range("K34").Value = 1054.5834
$ws.Range("M34").Value = -852.5834

$ws.Range("H58").Value = 1586.5227
$ws.Range("I58").Value = 1379.0312
$ws.Range("J58").Value = 2139.8333
$ws.Range("K58").Value = 1379.0312
$ws.Range("L58").Value = 2139.8333
$ws.Range("M58").Value = -1176.0312
$ws.Range("N58").Value = -2545.8333

$ws.Range("H132").Value = 1992.8379
$ws.Range("I132").Value = 1397.8846
$ws.Range("J132").Value = 3399.0908
$ws.Range("K132").Value = 4193.6538
$ws.Range("L132").Value = 10197.2724
$ws.Range("M132").Value = -1663.6538
$ws.Range("N132").Value = -15257.2724

$ws.Range("H134").Value = 1960.2963
$ws.Range("I134").Value = 1823.826
$ws.Range("J134").Value = 2745
$ws.Range("K134").Value = 5471.478
$ws.Range("L134").Value = 8235
$ws.Range("M134").Value = -2936.478
$ws.Range("N134").Value = -13305

$ws.Range("H136").Value = 1586.5227
$ws.Range("I136").Value = 1379.0312
$ws.Range("J136").Value = 2139.8333
$ws.Range("K136").Value = 4137.0936
$ws.Range("L136").Value = 6419.499899999999
$ws.Range("M136").Value = -1587.0936
$ws.Range("N136").Value = -11519.4999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1372019.2
$ws.Range("I131").Value = 4880
$ws.Range("J131").Value = 1614576.2
$ws.Range("K131").Value = 14640
$ws.Range("L131").Value = 4843728.6
$ws.Range("M131").Value = -9600
$ws.Range("N131").Value = -4853808.6

$ws.Range("H133").Value = 7573.759
$ws.Range("J133").Value = 8164.5415
$ws.Range("L133").Value = 24493.6245
$ws.Range("N133").Value = -34613.62450000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 27150
$ws.Range("I80").Value = 2200
$ws.Range("J80").Value = 35466.668
$ws.Range("K80").Value = 2200
$ws.Range("L80").Value = 35466.668
$ws.Range("M80").Value = -1202
$ws.Range("N80").Value = -37462.668

$ws.Range("H83").Value = 27150
$ws.Range("I83").Value = 2200
$ws.Range("J83").Value = 35466.668
$ws.Range("K83").Value = 11000
$ws.Range("L83").Value = 177333.34
$ws.Range("M83").Value = -6008
$ws.Range("N83").Value = -187317.34

$ws.Range("H102").Value = 1389.8948
$ws.Range("I102").Value = 1015.6667
$ws.Range("K102").Value = 1015.6667
$ws.Range("M102").Value = 606.3333

$ws.Range("H122").Value = 4411.6665
$ws.Range("I122").Value = 4746.136
$ws.Range("J122").Value = 2940
$ws.Range("K122").Value = 14238.408
$ws.Range("L122").Value = 8820
$ws.Range("M122").Value = -11788.408
$ws.Range("N122").Value = -13720

$ws.Range("H132").Value = 3433.3044
$ws.Range("I132").Value = 3257.375
$ws.Range("J132").Value = 3835.4285
$ws.Range("K132").Value = 9772.125
$ws.Range("L132").Value = 11506.2855
$ws.Range("M132").Value = -7242.125
$ws.Range("N132").Value = -16566.2855

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 7017.478
$ws.Range("I132").Value = 9700.333000000001
$ws.Range("J132").Value = 4090.7273
$ws.Range("K132").Value = 29100.999
$ws.Range("L132").Value = 12272.1819
$ws.Range("M132").Value = -26570.999
$ws.Range("N132").Value = -17332.1819

$ws.Range("H136").Value = 1849.25
$ws.Range("I136").Value = 819.1
$ws.Range("K136").Value = 2457.3
$ws.Range("M136").Value = 92.69999999999982

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 27000
$ws.Range("J92").Value = 27000
$ws.Range("L92").Value = 27000
$ws.Range("N92").Value = -31992

$ws.Range("H132").Value = 1368.4324
$ws.Range("I132").Value = 754.6070999999999
$ws.Range("J132").Value = 3278.111
$ws.Range("K132").Value = 2263.8213
$ws.Range("L132").Value = 9834.332999999999
$ws.Range("M132").Value = 266.1787000000004
$ws.Range("N132").Value = -14894.333

$ws.Range("H136").Value = 8846.115
$ws.Range("I136").Value = 9166.625
$ws.Range("K136").Value = 27499.875
$ws.Range("M136").Value = -24949.875
